$d = $word.ActiveDocument

$replacements = @(
    @("249÷6=", "226÷4="),
    @("771÷6=", "332÷4="),
    @("177÷4=", "353÷6="),
    @("661÷3=", "877÷8="),
    @("451÷5=", "225÷5="),
    @("606÷4=", "647÷2="),
    @("677÷8=", "230÷2="),
    @("442÷7=", "960÷2="),
    @("587÷9=", "550÷7="),
    @("820÷8=", "894÷2="),
    @("645÷5=", "374÷4="),
    @("258÷8=", "567÷7="),
    @("749÷5=", "759÷7="),
    @("344÷7=", "734÷6="),
    @("715÷3=", "119÷8="),
    @("774÷4=", "436÷8="),
    @("677÷6=", "619÷3="),
    @("808÷9=", "536÷6="),
    @("843÷8=", "874÷2="),
    @("745÷9=", "620÷2="),
    @("924÷8=", "536÷2="),
    @("923÷7=", "502÷6="),
    @("377÷3=", "256÷2="),
    @("281÷4=", "202÷9="),
    @("468÷4=", "573÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
